$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "44.472.09"
$ws.Range("E2").Value = "  +0.41%  "

# Row 3
$ws.Range("D3").Value = "2.221.21"
$ws.Range("E3").Value = "  -0.73%  "

# Row 4
$ws.Range("E4").Value = "  +0.25%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.10"
$ws.Range("E5").Value = "  -1.44%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.95"
$ws.Range("E6").Value = "  -4.32%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.559"
$ws.Range("E7").Value = "  -2.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  -4.69%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.73"
$ws.Range("E10").Value = "  -1.76%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0781"
$ws.Range("E11").Value = "  -3.21%  "

# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.95"
$ws.Range("E12").Value = "  -2.68%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  -0.33%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "2.362.55"
$ws.Range("E14").Value = "  +1.64%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.558.32"
$ws.Range("E15").Value = "  -0.72%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.806"
$ws.Range("E16").Value = "  -2.26%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.18"
$ws.Range("E17").Value = "  -1.76%  "

# Row 18
$ws.Range("D18").Value = "44.337.57"
$ws.Range("E18").Value = "  +0.88%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0911"
$ws.Range("E19").Value = "  -5.13%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.02"
$ws.Range("E20").Value = "  -5.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.33"
$ws.Range("E21").Value = "  -5.68%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "64.40"
$ws.Range("E22").Value = "  -1.65%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.24"
$ws.Range("E23").Value = "  -1.50%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.87"
$ws.Range("E24").Value = "  -7.47%  "

# Row 25
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("E26").Value = "  -4.35%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.27"
$ws.Range("E27").Value = "  +2.85%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.47"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.11"
$ws.Range("E29").Value = "  -8.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.55"
$ws.Range("E30").Value = "  -2.29%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.58"
$ws.Range("E31").Value = "  -4.35%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "145.99"
$ws.Range("E32").Value = "  -4.65%  "

# Row 33
$ws.Range("E33").Value = "  +0.57%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0758"
$ws.Range("E34").Value = "  -4.35%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.00"
$ws.Range("E35").Value = "  -2.32%  "

# Row 36
$ws.Range("E36").Value = "  -1.34%  "

# Row 37
$ws.Range("E37").Value = "  -3.41%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  +2.23%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.50"
$ws.Range("E39").Value = "  +2.82%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.23"
$ws.Range("E40").Value = "  -6.85%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.64"
$ws.Range("E41").Value = "  -3.78%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0287"
$ws.Range("E42").Value = "  -3.36%  "

# Row 43
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("D44").Value = "1.777.87"
$ws.Range("E44").Value = "  +3.32%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.72"
$ws.Range("E45").Value = "  +7.32%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "79.28"
$ws.Range("E46").Value = "  -3.71%  "

# Row 47
$ws.Range("E47").Value = "  -4.99%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "95.32"
$ws.Range("E48").Value = "  -3.90%  "

# Row 49
$ws.Range("E49").Value = "  -3.82%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.91"
$ws.Range("E50").Value = "  -0.09%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.31"
$ws.Range("E51").Value = "  -4.03%  "
